$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (15) with the SplitterWord setting used to separate
# generated DM messages, mirroring the existing Name/Value config rows.
$ws.Range("A15").Value = "SplitterWord"
$ws.Range("B15").Value = "I’m Visuals"

# Update the active cell selection to match the authored workbook state.
$ws.Range("H16").Select()
